# edit.ps1 - apply the "without design & red" changes described by the diff
#
# Summary of edits:
#  1. Paragraph starting "Через " gets an explicit en-US language on the
#     paragraph mark (w:pPr/w:rPr/w:lang).
#  2. Paragraph "Отсматривать динамику ученика в виде графика" is replaced
#     by a single run "Разобраться с непредвиденными ситуациями".
#  3. Paragraph "Telegram-бот дополнительно" is replaced by a single run
#     "Накатить дизайна".
#  4. Paragraph "Сделать аккаунт учителя ..." is replaced by a single run
#     "Добавить аналитики", and the _GoBack bookmark is moved so it comes
#     before the run instead of after it.

$d = $word.ActiveDocument
$origParaCount = $d.Paragraphs.Count

function Find-ParagraphByPrefix($doc, [string]$prefix) {
    $match = $null
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.StartsWith($prefix)) {
            $match = $p
        }
    }
    if ($null -eq $match) {
        throw "could not find a paragraph starting with '$prefix'"
    }
    return $match
}

function Wrap-Package([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# --- 1. "Через ..." paragraph: tag the paragraph mark as en-US -------------
$p1 = Find-ParagraphByPrefix $d "Через "
$frag1 = Wrap-Package (
    '<w:p><w:pPr><w:spacing w:after="960"/><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>' +
    '<w:r><w:t xml:space="preserve">Через </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>qt</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">заливка </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xlsx</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">файлов. Потом в </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>qt</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">зажимаешь кнопку и тебе выдает </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xlsx</w:t></w:r><w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>файл.</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
    '<w:r><w:t>Сделать класс ученика.</w:t></w:r>' +
    '</w:p>'
)
$p1.Range.InsertXML($frag1)

# --- 2. "Отсматривать динамику ..." -> "Разобраться с непредвиденными ситуациями"
$p2 = Find-ParagraphByPrefix $d "Отсматривать"
$frag2 = Wrap-Package (
    '<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Разобраться с непредвиденными ситуациями</w:t></w:r></w:p>'
)
$p2.Range.InsertXML($frag2)

# --- 3. "Telegram-бот дополнительно" -> "Накатить дизайна" -----------------
$p3 = Find-ParagraphByPrefix $d "Telegram-"
$frag3 = Wrap-Package (
    '<w:p><w:pPr><w:pStyle w:val="a3"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr></w:pPr>' +
    '<w:r><w:t>Накатить дизайна</w:t></w:r></w:p>'
)
$p3.Range.InsertXML($frag3)

# --- 4. "Сделать аккаунт ..." -> "Добавить аналитики" (+ bookmark moved) ---
# This is the last paragraph in the body; InsertXML-ing its whole Range
# (paragraph mark included) leaves a stray trailing empty paragraph because
# the body's final mark can't be consumed that way, so use Find/Replace for
# the text and the Bookmarks collection for the _GoBack relocation instead.
$p4 = Find-ParagraphByPrefix $d "Сделать аккаунт"
$p4Start = $p4.Range.Start
$d.Content.Find.Execute(
    "Сделать аккаунт учителя с историей, подпиской, и всем барахлом",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Добавить аналитики", 2) | Out-Null

# Move (or create) the _GoBack bookmark so it sits right at the start of
# the paragraph, i.e. before the run instead of after it.
$bmRange = $d.Range($p4Start, $p4Start)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- sanity checks -----------------------------------------------------
if ($d.Paragraphs.Count -ne $origParaCount) {
    throw "paragraph count changed: $origParaCount -> $($d.Paragraphs.Count)"
}
$bm = $d.Bookmarks.Item("_GoBack")
if ($bm.Start -ne $p4Start -or $bm.End -ne $p4Start) {
    throw "_GoBack bookmark is not at the expected position"
}
Write-Host "edit.ps1 completed: $($d.Paragraphs.Count) paragraphs, bookmark at $($bm.Start)"
